$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (modified: D,E,G,H,I,J,L,M,N,O,P,Q,R,S,T)
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 32.52516566666667
$ws.Cells.Item(2, 8).Value = 97.575497
$ws.Cells.Item(2, 9).Value = 0.06370711489344116
$ws.Cells.Item(2, 10).Value = 0.06370711489344116
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.009519999999999999
$ws.Cells.Item(2, 14).Value = 0.02856
$ws.Cells.Item(2, 15).Value = 0.0009583584527718872
$ws.Cells.Item(2, 16).Value = 0.0009583584527718872
$ws.Cells.Item(2, 17).Value = 0.3096395771466666
$ws.Cells.Item(2, 18).Value = 2.78675619432
$ws.Cells.Item(2, 19).Value = [double]"6.105425205983912E-05"
$ws.Cells.Item(2, 20).Value = [double]"6.105425205983912E-05"

# Row 3 (modified: D,E,G,H,I,J,K,M,N,O,P,Q,R,S,T)
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 32.52516566666667
$ws.Cells.Item(3, 8).Value = 97.575497
$ws.Cells.Item(3, 9).Value = 0.06370711489344116
$ws.Cells.Item(3, 10).Value = 0.06370711489344116
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 1.079146666666667
$ws.Cells.Item(3, 14).Value = 3.23744
$ws.Cells.Item(3, 15).Value = 0.1086354338004839
$ws.Cells.Item(3, 16).Value = 0.1086354338004838
$ws.Cells.Item(3, 17).Value = 35.09942411196445
$ws.Cells.Item(3, 18).Value = 315.89481700768
$ws.Cells.Item(3, 19).Value = 0.006920850062626247
$ws.Cells.Item(3, 20).Value = 0.006920850062626246

# Row 4 (modified: A,D,E,G,H,I,J,K,M,N,O,P,Q,R,S,T)
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 32.52516566666667
$ws.Cells.Item(4, 8).Value = 97.575497
$ws.Cells.Item(4, 9).Value = 0.06370711489344116
$ws.Cells.Item(4, 10).Value = 0.06370711489344116
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 8.844986
$ws.Cells.Item(4, 14).Value = 26.534958
$ws.Cells.Item(4, 15).Value = 0.8904062077467443
$ws.Cells.Item(4, 16).Value = 0.8904062077467442
$ws.Cells.Item(4, 17).Value = 287.6846349693474
$ws.Cells.Item(4, 18).Value = 2589.161714724126
$ws.Cells.Item(4, 19).Value = 0.05672521057875508
$ws.Cells.Item(4, 20).Value = 0.05672521057875507

# Row 5 (modified: D,E,G,H,I,J,L,M,N,O,P,Q,R,S,T)
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 12.155757
$ws.Cells.Item(5, 8).Value = 36.467271
$ws.Cells.Item(5, 9).Value = 0.02380950848190151
$ws.Cells.Item(5, 10).Value = 0.02380950848190151
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.009519999999999999
$ws.Cells.Item(5, 14).Value = 0.02856
$ws.Cells.Item(5, 15).Value = 0.0009583584527718872
$ws.Cells.Item(5, 16).Value = 0.0009583584527718872
$ws.Cells.Item(5, 17).Value = 0.11572280664
$ws.Cells.Item(5, 18).Value = 1.04150525976
$ws.Cells.Item(5, 19).Value = [double]"2.281804370997426E-05"
$ws.Cells.Item(5, 20).Value = [double]"2.281804370997426E-05"

# Row 6 (modified: A,E,G,H,I,J,K,M,N,O,P,Q,R,S,T)
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 12.155757
$ws.Cells.Item(6, 8).Value = 36.467271
$ws.Cells.Item(6, 9).Value = 0.02380950848190151
$ws.Cells.Item(6, 10).Value = 0.02380950848190151
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 1.079146666666667
$ws.Cells.Item(6, 14).Value = 3.23744
$ws.Cells.Item(6, 15).Value = 0.1086354338004839
$ws.Cells.Item(6, 16).Value = 0.1086354338004838
$ws.Cells.Item(6, 17).Value = 13.11784464736
$ws.Cells.Item(6, 18).Value = 118.06060182624
$ws.Cells.Item(6, 19).Value = 0.00258655628250767
$ws.Cells.Item(6, 20).Value = 0.00258655628250767

# Row 7 (modified: A,E,G,H,I,J,K,M,N,O,P,Q,R,S,T)
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 12.155757
$ws.Cells.Item(7, 8).Value = 36.467271
$ws.Cells.Item(7, 9).Value = 0.02380950848190151
$ws.Cells.Item(7, 10).Value = 0.02380950848190151
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 8.844986
$ws.Cells.Item(7, 14).Value = 26.534958
$ws.Cells.Item(7, 15).Value = 0.8904062077467443
$ws.Cells.Item(7, 16).Value = 0.8904062077467442
$ws.Cells.Item(7, 17).Value = 107.517500484402
$ws.Cells.Item(7, 18).Value = 967.6575043596182
$ws.Cells.Item(7, 19).Value = 0.02120013415568386
$ws.Cells.Item(7, 20).Value = 0.02120013415568386

# Row 8 (modified: A,D,E,G,H,I,J,L,M,N,O,P,Q,R,S,T)
$ws.Cells.Item(8, 1).Value = "M1"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 151.0671593333333
$ws.Cells.Item(8, 8).Value = 453.201478
$ws.Cells.Item(8, 9).Value = 0.2958955835892216
$ws.Cells.Item(8, 10).Value = 0.2958955835892216
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.009519999999999999
$ws.Cells.Item(8, 14).Value = 0.02856
$ws.Cells.Item(8, 15).Value = 0.0009583584527718872
$ws.Cells.Item(8, 16).Value = 0.0009583584527718872
$ws.Cells.Item(8, 17).Value = 1.438159356853333
$ws.Cells.Item(8, 18).Value = 12.94343421168
$ws.Cells.Item(8, 19).Value = 0.0002835740336706011
$ws.Cells.Item(8, 20).Value = 0.0002835740336706011

# Row 9 (modified: A,D,E,G,H,I,J,K,M,N,O,P,Q,R,S,T)
$ws.Cells.Item(9, 1).Value = "M1"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 151.0671593333333
$ws.Cells.Item(9, 8).Value = 453.201478
$ws.Cells.Item(9, 9).Value = 0.2958955835892216
$ws.Cells.Item(9, 10).Value = 0.2958955835892216
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 1.079146666666667
$ws.Cells.Item(9, 14).Value = 3.23744
$ws.Cells.Item(9, 15).Value = 0.1086354338004839
$ws.Cells.Item(9, 16).Value = 0.1086354338004838
$ws.Cells.Item(9, 17).Value = 163.0236214373689
$ws.Cells.Item(9, 18).Value = 1467.21259293632
$ws.Cells.Item(9, 19).Value = 0.03214474508286243
$ws.Cells.Item(9, 20).Value = 0.03214474508286242

# Row 10 (modified: A,D,E,G,H,I,J,K,M,N,O,P,Q,R,S,T)
$ws.Cells.Item(10, 1).Value = "M1"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 151.0671593333333
$ws.Cells.Item(10, 8).Value = 453.201478
$ws.Cells.Item(10, 9).Value = 0.2958955835892216
$ws.Cells.Item(10, 10).Value = 0.2958955835892216
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 8.844986
$ws.Cells.Item(10, 14).Value = 26.534958
$ws.Cells.Item(10, 15).Value = 0.8904062077467443
$ws.Cells.Item(10, 16).Value = 0.8904062077467442
$ws.Cells.Item(10, 17).Value = 1336.186909363103
$ws.Cells.Item(10, 18).Value = 12025.68218426793
$ws.Cells.Item(10, 19).Value = 0.2634672644726886
$ws.Cells.Item(10, 20).Value = 0.2634672644726886

# Row 11 (modified: A,D,E,G,H,I,J,L,M,N,O,P,Q,R,S,T)
$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 162.458898
$ws.Cells.Item(11, 8).Value = 487.376694
$ws.Cells.Item(11, 9).Value = 0.3182086076491469
$ws.Cells.Item(11, 10).Value = 0.318208607649147
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.009519999999999999
$ws.Cells.Item(11, 14).Value = 0.02856
$ws.Cells.Item(11, 15).Value = 0.0009583584527718872
$ws.Cells.Item(11, 16).Value = 0.0009583584527718872
$ws.Cells.Item(11, 17).Value = 1.54660870896
$ws.Cells.Item(11, 18).Value = 13.91947838064
$ws.Cells.Item(11, 19).Value = 0.000304957908885333
$ws.Cells.Item(11, 20).Value = 0.000304957908885333

# Row 12 (modified: A,E,G,H,I,J,K,M,N,O,P,Q,R,S,T)
$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 162.458898
$ws.Cells.Item(12, 8).Value = 487.376694
$ws.Cells.Item(12, 9).Value = 0.3182086076491469
$ws.Cells.Item(12, 10).Value = 0.318208607649147
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 1.079146666666667
$ws.Cells.Item(12, 14).Value = 3.23744
$ws.Cells.Item(12, 15).Value = 0.1086354338004839
$ws.Cells.Item(12, 16).Value = 0.1086354338004838
$ws.Cells.Item(12, 17).Value = 175.31697824704
$ws.Cells.Item(12, 18).Value = 1577.85280422336
$ws.Cells.Item(12, 19).Value = 0.03456873013101305
$ws.Cells.Item(12, 20).Value = 0.03456873013101305

# Row 13 (modified: A,E,G,H,I,J,K,M,N,O,P,Q,R,S,T)
$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 162.458898
$ws.Cells.Item(13, 8).Value = 487.376694
$ws.Cells.Item(13, 9).Value = 0.3182086076491469
$ws.Cells.Item(13, 10).Value = 0.318208607649147
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 8.844986
$ws.Cells.Item(13, 14).Value = 26.534958
$ws.Cells.Item(13, 15).Value = 0.8904062077467443
$ws.Cells.Item(13, 16).Value = 0.8904062077467442
$ws.Cells.Item(13, 17).Value = 1436.946678385428
$ws.Cells.Item(13, 18).Value = 12932.52010546886
$ws.Cells.Item(13, 19).Value = 0.2833349196092486
$ws.Cells.Item(13, 20).Value = 0.2833349196092486

# Row 14 (new)
$ws.Cells.Item(14, 1).Value = "Neutro"
$ws.Cells.Item(14, 2).Value = "Tgfb1"
$ws.Cells.Item(14, 3).Value = "Itgb8"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 142.5917893333333
$ws.Cells.Item(14, 8).Value = 427.775368
$ws.Cells.Item(14, 9).Value = 0.2792948573734662
$ws.Cells.Item(14, 10).Value = 0.2792948573734662
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.009519999999999999
$ws.Cells.Item(14, 14).Value = 0.02856
$ws.Cells.Item(14, 15).Value = 0.0009583584527718872
$ws.Cells.Item(14, 16).Value = 0.0009583584527718872
$ws.Cells.Item(14, 17).Value = 1.357473834453333
$ws.Cells.Item(14, 18).Value = 12.21726451008
$ws.Cells.Item(14, 19).Value = 0.00026766458737958
$ws.Cells.Item(14, 20).Value = 0.00026766458737958

# Row 15 (new)
$ws.Cells.Item(15, 1).Value = "Neutro"
$ws.Cells.Item(15, 2).Value = "Tgfb1"
$ws.Cells.Item(15, 3).Value = "Itgb8"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 142.5917893333333
$ws.Cells.Item(15, 8).Value = 427.775368
$ws.Cells.Item(15, 9).Value = 0.2792948573734662
$ws.Cells.Item(15, 10).Value = 0.2792948573734662
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 1.079146666666667
$ws.Cells.Item(15, 14).Value = 3.23744
$ws.Cells.Item(15, 15).Value = 0.1086354338004839
$ws.Cells.Item(15, 16).Value = 0.1086354338004838
$ws.Cells.Item(15, 17).Value = 153.8774541531022
$ws.Cells.Item(15, 18).Value = 1384.89708737792
$ws.Cells.Item(15, 19).Value = 0.03034131798901077
$ws.Cells.Item(15, 20).Value = 0.03034131798901076

# Row 16 (new)
$ws.Cells.Item(16, 1).Value = "Neutro"
$ws.Cells.Item(16, 2).Value = "Tgfb1"
$ws.Cells.Item(16, 3).Value = "Itgb8"
$ws.Cells.Item(16, 4).Value = "sCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 142.5917893333333
$ws.Cells.Item(16, 8).Value = 427.775368
$ws.Cells.Item(16, 9).Value = 0.2792948573734662
$ws.Cells.Item(16, 10).Value = 0.2792948573734662
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 8.844986
$ws.Cells.Item(16, 14).Value = 26.534958
$ws.Cells.Item(16, 15).Value = 0.8904062077467443
$ws.Cells.Item(16, 16).Value = 0.8904062077467442
$ws.Cells.Item(16, 17).Value = 1261.222380368283
$ws.Cells.Item(16, 18).Value = 11351.00142331455
$ws.Cells.Item(16, 19).Value = 0.2486858747970758
$ws.Cells.Item(16, 20).Value = 0.2486858747970758

# Row 17 (new)
$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Tgfb1"
$ws.Cells.Item(17, 3).Value = "Itgb8"
$ws.Cells.Item(17, 4).Value = "ECs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 9.743353333333333
$ws.Cells.Item(17, 8).Value = 29.23006
$ws.Cells.Item(17, 9).Value = 0.01908432801282251
$ws.Cells.Item(17, 10).Value = 0.01908432801282251
$ws.Cells.Item(17, 11).Value = 1
$ws.Cells.Item(17, 12).Value = 0.3333333333333333
$ws.Cells.Item(17, 13).Value = 0.009519999999999999
$ws.Cells.Item(17, 14).Value = 0.02856
$ws.Cells.Item(17, 15).Value = 0.0009583584527718872
$ws.Cells.Item(17, 16).Value = 0.0009583584527718872
$ws.Cells.Item(17, 17).Value = 0.09275672373333332
$ws.Cells.Item(17, 18).Value = 0.8348105136
$ws.Cells.Item(17, 19).Value = [double]"1.828962706655976E-05"
$ws.Cells.Item(17, 20).Value = [double]"1.828962706655977E-05"

# Row 18 (new)
$ws.Cells.Item(18, 1).Value = "sCs"
$ws.Cells.Item(18, 2).Value = "Tgfb1"
$ws.Cells.Item(18, 3).Value = "Itgb8"
$ws.Cells.Item(18, 4).Value = "FAPs"
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 9.743353333333333
$ws.Cells.Item(18, 8).Value = 29.23006
$ws.Cells.Item(18, 9).Value = 0.01908432801282251
$ws.Cells.Item(18, 10).Value = 0.01908432801282251
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 12).Value = 1
$ws.Cells.Item(18, 13).Value = 1.079146666666667
$ws.Cells.Item(18, 14).Value = 3.23744
$ws.Cells.Item(18, 15).Value = 0.1086354338004839
$ws.Cells.Item(18, 16).Value = 0.1086354338004838
$ws.Cells.Item(18, 17).Value = 10.51450727182222
$ws.Cells.Item(18, 18).Value = 94.63056544640001
$ws.Cells.Item(18, 19).Value = 0.002073234252463699
$ws.Cells.Item(18, 20).Value = 0.002073234252463699

# Row 19 (new)
$ws.Cells.Item(19, 1).Value = "sCs"
$ws.Cells.Item(19, 2).Value = "Tgfb1"
$ws.Cells.Item(19, 3).Value = "Itgb8"
$ws.Cells.Item(19, 4).Value = "sCs"
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 9.743353333333333
$ws.Cells.Item(19, 8).Value = 29.23006
$ws.Cells.Item(19, 9).Value = 0.01908432801282251
$ws.Cells.Item(19, 10).Value = 0.01908432801282251
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 8.844986
$ws.Cells.Item(19, 14).Value = 26.534958
$ws.Cells.Item(19, 15).Value = 0.8904062077467443
$ws.Cells.Item(19, 16).Value = 0.8904062077467442
$ws.Cells.Item(19, 17).Value = 86.17982382638667
$ws.Cells.Item(19, 18).Value = 775.6184144374802
$ws.Cells.Item(19, 19).Value = 0.01699280413329225
$ws.Cells.Item(19, 20).Value = 0.01699280413329225

